# Actualización desde MV -datos-
# Adds a new "07-09-2021" series row (row 19) to the daily-rate table,
# mirroring the shape of the existing rows that only carry Serie / Cupo /
# Total monto adjudicado (e.g. row 18, "02-09-2021").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 19
$lastRow = $newRow - 1

# Force column A to be stored as text so the date-like label "07-09-2021"
# becomes a shared string (matching every other "Serie" cell in the sheet)
# instead of being auto-converted into a date serial number.
$ws.Range("A" + $newRow).NumberFormat = "@"
$ws.Range("A" + $newRow).Value = "07-09-2021"

# Re-apply the same (unformatted/default) style used by the rest of the
# data rows, so the new cell doesn't keep a stray text-number-format style.
$ws.Range("A" + $newRow).Style = $ws.Range("A" + $lastRow).Style

$ws.Range("B" + $newRow).Value = 10000
$ws.Range("D" + $newRow).Value = 0
